$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column D ("description") between "name" and "path" ---
$ws.Columns.Item(4).Insert()

# Header
$ws.Cells.Item(1,4).Value = "description"

# Description values (row 2 uses a distinct "pdf wiki" label; the rest mirror the "name" column)
$ws.Cells.Item(2,4).Value = "pdf wiki"
$ws.Cells.Item(3,4).Value = "Tourisme exemple"
$ws.Cells.Item(4,4).Value = "pdf online"
$ws.Cells.Item(5,4).Value = "BEVNAT: Fiche signalétique"
$ws.Cells.Item(6,4).Value = "STATPOP: Fiche signalétique"
$ws.Cells.Item(7,4).Value = "BEVNAT: Liste des variables"
$ws.Cells.Item(8,4).Value = "Communiqué de presse population"

# --- Expand the table (ListObject) to include the new column ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F8"))

# Re-assert the header text for the columns after the insertion point so the
# table's column metadata re-binds to the right header cells/names.
$ws.Cells.Item(1,5).Value = "path"
$ws.Cells.Item(1,6).Value = "last_update"

# --- Column width for the new "description" column (old "path"/"last_update"
# columns keep their own widths automatically, shifted right by the insert) ---
$ws.Columns.Item(4).ColumnWidth = 21.1666666666667

# --- Selection moves to D3 ---
$ws.Range("D3").Select()
